# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 124
$ws1.Range("F3").Value  = 2143
$ws1.Range("F5").Value  = 11197
$ws1.Range("F9").Value  = 215
$ws1.Range("F10").Value = 11101
$ws1.Range("F11").Value = 447
$ws1.Range("F15").Value = 5562
$ws1.Range("F17").Value = 3441
$ws1.Range("F18").Value = 6

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 124
$ws4.Range("F3").Value  = 2143
$ws4.Range("F7").Value  = 11197
$ws4.Range("F11").Value = 215
$ws4.Range("F12").Value = 11101
$ws4.Range("F13").Value = 447
$ws4.Range("F17").Value = 5562
$ws4.Range("F19").Value = 3441
$ws4.Range("F20").Value = 6
